$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append new weekly observation row (2023-07-26 / 543.597) ---
$dataSheet = $wb.Worksheets.Item("Data")

# Copy the formatting of the previous date cell (A95) onto the new date
# cell (A96) before writing values, so the new row matches the existing
# date-column style (centered, bordered, yyyy-mm-dd date format).
$dataSheet.Cells.Item(95, 1).Copy()
$dataSheet.Cells.Item(96, 1).PasteSpecial(-4122)   # xlPasteFormats

$dataSheet.Cells.Item(96, 1).Value = 45133
$dataSheet.Cells.Item(96, 2).Value = 543.597

# --- Sheet "SeriesInfo": refresh the metadata pulled from the FRED API ---
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# B3/B4/B7 hold plain textual dates (no time component) which Excel would
# otherwise auto-convert to date serials on assignment. Force the cell to
# text first, write the string, then reset the formatting back to the
# sheet's default (matching an untouched text cell) so no stray
# number-format is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $infoSheet.Range("B2").Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

Set-TextValue $infoSheet.Range("B3") "2023-08-03"
Set-TextValue $infoSheet.Range("B4") "2023-08-03"
Set-TextValue $infoSheet.Range("B7") "2023-07-26"

# B14 already contains a UTC-offset timestamp string that Excel leaves
# alone, so a direct assignment is sufficient.
$infoSheet.Range("B14").Value = "2023-07-27 15:34:02-05"

# B15 is a genuine number.
$infoSheet.Range("B15").Value = 82
